$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '290.78'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-8.12%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '40.40'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-1.64%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.030'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-3.06%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07292'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-4.43%'
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.285'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.82%'
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.563'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-5.38%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9205'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-1.28%'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1160'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-6.88%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1729'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-5.33%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08662'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-4.52%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04181'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '1.31%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1053'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.25%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001275'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.21%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005890'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.51%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.403'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1.50%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.358'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-2.76%'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-2.52%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-6.28%'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '1.62%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.35%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.03858'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-4.52%'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.31%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.003794'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-6.58%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.32%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0003727'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-95.02%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02309'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-6.44%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.04946'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-5.37%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.006538'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '202.22%'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.47%'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-1.68%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.007360'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '3.90%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007073'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-14.05%'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-15.35%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006413'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-4.04%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.45%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.02986'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-91.18%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.55%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.45%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002003'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.45%'
